# Adds two new "Title and Content" slides to the end of the deck:
#   9  -> "Test Case Document"
#   10 -> "RTM"
# matching the commit "updated ppt to include RTM".

$p = $ppt.ActivePresentation

# --- Slide 9: Test Case Document ---------------------------------------
$s9 = $p.Slides.Add(9, 2)   # ppLayoutText -> same "Title and Content" layout as the other slides

$s9.Shapes.Item(1).TextFrame.TextRange.Text = "Test Case Document"

$s9Body = $s9.Shapes.Item(2).TextFrame.TextRange
$s9Body.Text = "Outlines all of the specific test cases for all user stories"
[void]$s9Body.InsertAfter("`rUseful for the RTM")
[void]$s9Body.InsertAfter("`rYou don" + [char]0x2019 + "t have to go into code base to have knowledge of what tests are exiting and if they pass or not, etc.")
[void]$s9Body.InsertAfter("`rhttps://www.guru99.com/download-sample-test-case-template-with-explanation-of-important-fields.html")

# --- Slide 10: RTM -------------------------------------------------------
$s10 = $p.Slides.Add(10, 2)

$s10.Shapes.Item(1).TextFrame.TextRange.Text = "RTM"

$s10Body = $s10.Shapes.Item(2).TextFrame.TextRange
$s10Body.Text = "Requirements Traceability Matrix outlines all test cases and compares them to the desired requirements"
[void]$s10Body.InsertAfter("`rIt marks if they have passed or failed their test")
[void]$s10Body.InsertAfter("`rSingle location where you can show that you succeeded at your job as tester")
[void]$s10Body.InsertAfter("`rhttps://www.guru99.com/traceability-matrix.html")
